$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F values for a handful of rows (simulating a different
# random missing-data pattern for this re-run of the imputation script).
$ws.Range("F3").Value = 17.64
$ws.Range("F5").Value = ""
$ws.Range("F21").Value = 16.58
$ws.Range("F23").Value = ""

# Remove the "RM 232" row entirely (row 26); rows below shift up.
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" row entirely (now row 27 after the previous delete);
# rows below shift up again.
$ws.Rows.Item(27).Delete()

# The row that used to be "SC 193" (originally row 34) is now row 32;
# give it a value for column F.
$ws.Range("F32").Value = 17.39
